# Update the default category encoder to include a delimiter.
# The "Categories" column (F) value for the sample row gets a second,
# comma-delimited category appended.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the extra category ("Biochemistry") to the existing
# "Computer Software" value in the Categories column (F2).
$ws.Range("F2").Value = "Computer Software, Biochemistry"

# Leave the freshly-edited cell selected/active, as it would be right
# after a user types the new value in and moves off of it.
$ws.Range("F2").Select()
